# Apply FlashScore odds/score updates for 2025-04-22 workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 6.5
$ws.Range("H3").Value = 5.5
$ws.Range("R3").Value = 1.62
$ws.Range("S3").Value = 2.2
$ws.Range("T3").Value = 23
$ws.Range("Z3").Value = 21
$ws.Range("AF3").Value = 9
$ws.Range("AG3").Value = 9
# Row 4
$ws.Range("N4").Value = 2.1
$ws.Range("O4").Value = 1.73
# Row 5
$ws.Range("K5").Value = 19
$ws.Range("W5").Value = 7.5
$ws.Range("Y5").Value = 26
$ws.Range("Z5").Value = 19
# Row 9
$ws.Range("AA9").Value = 6
$ws.Range("AB9").Value = 23
$ws.Range("AC9").Value = 101
$ws.Range("AE9").Value = 7
$ws.Range("AF9").Value = 17
# Row 11
$ws.Range("K11").Value = 23
# Row 16
$ws.Range("G16").Value = 1.9
$ws.Range("I16").Value = 3.75
$ws.Range("AG16").Value = 13
# Row 18
$ws.Range("P18").Value = 1.53
$ws.Range("Q18").Value = 2.38
$ws.Range("R18").Value = 2.1
$ws.Range("S18").Value = 1.67
$ws.Range("AJ18").Value = 51
# Row 20
$ws.Range("G20").Value = 3.3
$ws.Range("H20").Value = 3.1
$ws.Range("I20").Value = 2.15
$ws.Range("L20").Value = 1.34
$ws.Range("M20").Value = 2.75
$ws.Range("N20").Value = 2
$ws.Range("O20").Value = 1.65
$ws.Range("T20").Value = 9.5
$ws.Range("U20").Value = 17.5
$ws.Range("V20").Value = 11.5
$ws.Range("Y20").Value = 37
$ws.Range("AC20").Value = 70
$ws.Range("AE20").Value = 6.8
$ws.Range("AF20").Value = 10
$ws.Range("AG20").Value = 8.75
$ws.Range("AH20").Value = 21
$ws.Range("AJ20").Value = 32
# Row 21
$ws.Range("G21").Value = 2.37
$ws.Range("I21").Value = 2.9
$ws.Range("L21").Value = 1.34
$ws.Range("M21").Value = 2.75
$ws.Range("N21").Value = 2
$ws.Range("O21").Value = 1.65
$ws.Range("P21").Value = 1.4
$ws.Range("Q21").Value = 2.5
$ws.Range("R21").Value = 1.75
$ws.Range("S21").Value = 1.85
$ws.Range("T21").Value = 7.4
$ws.Range("V21").Value = 9.25
$ws.Range("W21").Value = 25
$ws.Range("X21").Value = 20
$ws.Range("Y21").Value = 32
$ws.Range("Z21").Value = 8.5
$ws.Range("AB21").Value = 14
$ws.Range("AC21").Value = 70
$ws.Range("AD21").Value = 600
$ws.Range("AF21").Value = 14.5
$ws.Range("AG21").Value = 10.5
$ws.Range("AH21").Value = 37
$ws.Range("AI21").Value = 27
$ws.Range("AJ21").Value = 37
# Row 22
$ws.Range("G22").Value = 1.39
$ws.Range("H22").Value = 3.65
$ws.Range("I22").Value = 10
$ws.Range("T22").Value = 5.4
$ws.Range("AB22").Value = 21
$ws.Range("AE22").Value = 21
$ws.Range("AF22").Value = 80
$ws.Range("AH22").Value = 400
$ws.Range("AI22").Value = 175
# Row 31
$ws.Range("P31").Value = 1.25
$ws.Range("Q31").Value = 3.75
$ws.Range("R31").Value = 1.44
$ws.Range("S31").Value = 2.63
$ws.Range("T31").Value = 13
$ws.Range("V31").Value = 10
$ws.Range("W31").Value = 26
$ws.Range("AB31").Value = 11
$ws.Range("AF31").Value = 17
# Row 36
$ws.Range("G36").Value = 2.25
$ws.Range("H36").Value = 3.3
$ws.Range("I36").Value = 2.77
$ws.Range("N36").Value = 1.62
$ws.Range("O36").Value = 2.02
$ws.Range("R36").Value = 1.53
$ws.Range("S36").Value = 2.35
$ws.Range("U36").Value = 10.75
$ws.Range("V36").Value = 7.6
$ws.Range("W36").Value = 19
$ws.Range("X36").Value = 14
$ws.Range("Y36").Value = 17.5
$ws.Range("Z36").Value = 12
$ws.Range("AA36").Value = 5.9
$ws.Range("AE36").Value = 9.5
$ws.Range("AF36").Value = 14
$ws.Range("AG36").Value = 8.5
$ws.Range("AH36").Value = 27
$ws.Range("AI36").Value = 17
$ws.Range("AJ36").Value = 19.5
# Row 45
$ws.Range("I45").Value = 1.35
$ws.Range("K45").Value = 10.5
$ws.Range("L45").Value = 1.13
$ws.Range("M45").Value = 4.5
$ws.Range("N45").Value = 1.53
$ws.Range("O45").Value = 2.38
# Row 46
$ws.Range("I46").Value = 1.37
$ws.Range("J46").Value = 1.01
$ws.Range("L46").Value = 1.13
# Row 50
$ws.Range("G50").Value = 1.44
$ws.Range("H50").Value = 4.4
$ws.Range("I50").Value = 6.3
$ws.Range("P50").Value = 1.33
$ws.Range("Q50").Value = 3.05
$ws.Range("R50").Value = 1.85
$ws.Range("S50").Value = 1.85
$ws.Range("T50").Value = 7.4
$ws.Range("U50").Value = 7
$ws.Range("W50").Value = 9.75
$ws.Range("X50").Value = 11.25
$ws.Range("Y50").Value = 25
$ws.Range("AB50").Value = 18
$ws.Range("AC50").Value = 80
$ws.Range("AD50").Value = 600
$ws.Range("AE50").Value = 17.5
$ws.Range("AF50").Value = 40
$ws.Range("AG50").Value = 19.5
$ws.Range("AH50").Value = 120
$ws.Range("AI50").Value = 65
$ws.Range("AJ50").Value = 60
# Row 51
$ws.Range("G51").Value = 8.25
$ws.Range("H51").Value = 4.7
$ws.Range("I51").Value = 1.34
$ws.Range("K51").Value = 9
$ws.Range("L51").Value = 1.17
$ws.Range("M51").Value = 4.45
$ws.Range("N51").Value = 1.52
$ws.Range("O51").Value = 2.35
$ws.Range("P51").Value = 1.28
$ws.Range("Q51").Value = 3.3
$ws.Range("R51").Value = 1.8
$ws.Range("S51").Value = 1.91
$ws.Range("T51").Value = 26
$ws.Range("U51").Value = 65
$ws.Range("V51").Value = 25
$ws.Range("W51").Value = 200
$ws.Range("X51").Value = 90
$ws.Range("Y51").Value = 65
$ws.Range("Z51").Value = 9
$ws.Range("AA51").Value = 9.5
$ws.Range("AB51").Value = 18
$ws.Range("AC51").Value = 70
$ws.Range("AD51").Value = 450
$ws.Range("AE51").Value = 8.25
$ws.Range("AF51").Value = 7.1
$ws.Range("AG51").Value = 8.25
$ws.Range("AH51").Value = 9
$ws.Range("AI51").Value = 10.5
# Row 52
$ws.Range("G52").Value = 4.5
$ws.Range("H52").Value = 4.05
$ws.Range("I52").Value = 1.65
$ws.Range("M52").Value = 4.25
$ws.Range("O52").Value = 2.3
$ws.Range("P52").Value = 1.3
$ws.Range("Q52").Value = 3.25
$ws.Range("R52").Value = 1.6
$ws.Range("S52").Value = 2.2
$ws.Range("T52").Value = 16
$ws.Range("U52").Value = 28
$ws.Range("V52").Value = 14.5
$ws.Range("W52").Value = 75
$ws.Range("X52").Value = 37
$ws.Range("Y52").Value = 37
$ws.Range("AA52").Value = 8.25
$ws.Range("AB52").Value = 14
$ws.Range("AF52").Value = 9.25
$ws.Range("AH52").Value = 13.5
$ws.Range("AI52").Value = 11.75
$ws.Range("AJ52").Value = 19.5
